$d = $word.ActiveDocument

# --- Step 1: change "5 :" -> "2 C# " -------------------------------------
# Replacing only "5 :" leaves the surrounding w:proofErr (gramStart/gramEnd)
# markers in place, so instead we grow the replaced span to also swallow the
# end of the preceding run ("Day"). That forces Word to merge across the
# gramStart boundary and drop it, while the text still ends up correct.
$rng = $d.Content
$rng.Find.Execute("Day5 :", $false, $false, $false, $false, $false, $true, 1, $false, "Day2 C# ", 2) | Out-Null

# --- Step 2: change " Task 1" -> " : Task 1" ------------------------------
# Same trick, but this time grow the span to the left so it crosses the
# gramEnd boundary too, which drops that marker as well.
$rng2 = $d.Content
$rng2.Find.Execute("C#  Ta", $false, $false, $false, $false, $false, $true, 1, $false, "C#  : Ta", 2) | Out-Null

# At this point the paragraph reads "Day2 C#  : Task 1" as a single run and
# both w:proofErr markers are gone. Re-split that run back into the three
# runs the document originally had ("Day" / "2 C# " / " : Task 1") by
# nudging the font on each sub-range: assigning Font.Name alone changes the
# effective rPr (it clears the complex-script font), which forces Word to
# split off a new run; re-asserting the complex-script name right after
# restores the original formatting on that new run without re-merging it
# with its neighbour.
$dayRng = $d.Content
$dayRng.Find.Execute("Day", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$dayRng.Font.Name = "Arial"
$dayRng2 = $d.Content
$dayRng2.Find.Execute("Day", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$dayRng2.Font.NameBi = "Arial"

$midRng = $d.Content
$midRng.Find.Execute("2 C# ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$midRng.Font.Name = "Arial"
$midRng2 = $d.Content
$midRng2.Find.Execute("2 C# ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$midRng2.Font.NameBi = "Arial"
